$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.169.42"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").Value = "3.511.93"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("D7").Value = "3.512.62"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("E11").Value = "  +6.29%  "

$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "4.095.36"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "3.507.07"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "67.270.33"
$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.87"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  -2.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "436.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.28%  "

$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("D25").Value = "3.652.70"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  -4.58%  "

$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  -4.71%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("E33").Value = "  -2.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.59"
$ws.Range("D34").Style = "Normal"

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.02%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.12%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "176.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0903"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("E43").Value = "  -10.64%  "

$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("E45").Value = "  -1.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.42%  "

$ws.Range("E47").Value = "  -5.88%  "

$ws.Range("E48").Value = "  -2.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("E51").Value = "  -2.45%  "
